$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

# --- Sheet "About" ---
$wsAbout = $wb.Worksheets.Item("About")

$a2 = $wsAbout.Range("A2")
$a2.Value = $a2.Value().Replace($oldVersion, $newVersion)

$a6 = $wsAbout.Range("A6")
$a6.Value = $a6.Value().Replace($oldVersion, $newVersion)

# --- Sheet "Boundaries and methane sources" ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

$lastRow = $wsData.Cells.Item($wsData.Rows.Count, 19).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsData.Cells.Item($r, 19)
    $v = $cell.Value()
    if ($v -ne $null -and $v -ne "") {
        $cell.Value = $v.Replace($oldVersion, $newVersion)
    }
}
